# Update "Shot Classification Models.xlsx" - add the EfficientNetB7
# (3,15,224,224,3) training run recorded in row 10, and refresh the
# sheet view (zoom + selection) to match the saved state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New data row (row 10) ---------------------------------------------
$ws.Range("B10").Value2 = 45252
$ws.Range("B10").NumberFormat = $ws.Range("B8").NumberFormat

$ws.Range("C10").Value = "EfficientNetB7"
$ws.Range("D10").Value = "Yes"
$ws.Range("E10").Value = "No"
$ws.Range("F10").Value = "(3, 15, 224, 224, 3)"

$ws.Range("I8").Copy() | Out-Null
$ws.Range("I10").PasteSpecial(-4122) | Out-Null
$ws.Range("I10").Value = "NA"

$ws.Range("J10").Value = "approx 3 hr"
$ws.Range("L10").Value = 0.705
$ws.Range("M10").Value = 0.495

# --- Column width refresh (auto-fit picked up the new row 10 content) --
$ws.Columns.Item(3).ColumnWidth = 17.666666666666668
$ws.Columns.Item(4).ColumnWidth = 17.333333333333332
$ws.Columns.Item(5).ColumnWidth = 21.833333333333332
$ws.Columns.Item(6).ColumnWidth = 15.666666666666666
$ws.Columns.Item(9).ColumnWidth = 17.5
$ws.Columns.Item(10).ColumnWidth = 17.166666666666668
$ws.Columns.Item(12).ColumnWidth = 13.333333333333334
$ws.Columns.Item(13).ColumnWidth = 12.5

# --- Sheet view: zoom + selected cell -----------------------------------
$excel.ActiveWindow.Zoom = 67
$ws.Range("M11").Select() | Out-Null
